$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for team record columns
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Match the header formatting used by the rest of row 1 (bold, bordered,
# centered / top-aligned)
$hdr = $ws.Range("AC1:AE1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

# Team record values for every player row (2-43)
$ws.Range("AC2:AC43").Value = 98
$ws.Range("AD2:AD43").Value = 64
$ws.Range("AE2:AE43").Value = 0

Write-Output "done"
